$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Weekly price-sheet update ---
# Two new price records were added for "Ramas de apio" (Vega Modelo de Temuco).
# The first is inserted right after the existing first record (becomes row 3,
# pushing the rest of the table down by one row). The second is appended
# after the (now shifted) last historical record (becomes row 17, pushing the
# former last row down to row 18).

# 1) Insert a new row at row 3 and populate it.
$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = 10
$ws.Range("B3").Value = 'Vega Modelo de Temuco'
$ws.Range("C3").Value = 'La Araucanía'
$ws.Range("D3").Value = 44649
$ws.Range("E3").Value = 9
$ws.Range("F3").Value = 100112017
$ws.Range("G3").Value = 'Ramas de apio'
$ws.Range("H3").Value = 'Sin especificar'
$ws.Range("I3").Value = 'Primera'
$ws.Range("J3").Value = 20
$ws.Range("K3").Value = 5000
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = 5000
$ws.Range("N3").Value = '$/paquete'
$ws.Range("O3").Value = 'Región de La Araucanía'
$ws.Range("P3").Value = 5000
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = 'Hortaliza'

# 2) Insert a new row at row 17 (after the shift above) and populate it.
$ws.Rows.Item(17).Insert()

$ws.Range("A17").Value = 10
$ws.Range("B17").Value = 'Vega Modelo de Temuco'
$ws.Range("C17").Value = 'La Araucanía'
$ws.Range("D17").Value = 44390
$ws.Range("E17").Value = 9
$ws.Range("F17").Value = 100112017
$ws.Range("G17").Value = 'Ramas de apio'
$ws.Range("H17").Value = 'Sin especificar'
$ws.Range("I17").Value = 'Primera'
$ws.Range("J17").Value = 55
$ws.Range("K17").Value = 6000
$ws.Range("L17").Value = 6000
$ws.Range("M17").Value = 6000
$ws.Range("N17").Value = '$/paquete'
$ws.Range("O17").Value = 'Región de La Araucanía'
$ws.Range("P17").Value = 6000
$ws.Range("Q17").Value = 1
$ws.Range("R17").Value = 'Hortaliza'
